$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values (M, SD, Median, Maximum columns) with new data
$ws.Range("B2").Value = 2.35554367269041
$ws.Range("C2").Value = 3.5913753222137799
$ws.Range("E2").Value = 1.575

$ws.Range("B3").Value = 1.5718970429633099
$ws.Range("C3").Value = 1.16991710716116
$ws.Range("E3").Value = 1.272

$ws.Range("B4").Value = 3.47158829215897
$ws.Range("C4").Value = 3.4771365855502898

$ws.Range("B5").Value = 1.94247102315491
$ws.Range("C5").Value = 1.71678037018383
$ws.Range("E5").Value = 1.5249999999999999

$ws.Range("B6").Value = 2.0231033901472002
$ws.Range("C6").Value = 1.60335224435627
$ws.Range("F6").Value = 19.068999999999999

# Update selection to match new active cell/selection
$ws.Range("C19").Select()
